$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Insert three new rows right before current row 61 (shifts old rows 61-180 down to 64-183)
$ws.Rows("61:63").Insert()

# New indicator keys to insert
$newKeys = @("IND_12_CHANGE", "IND_12_DELTA_DT", "IND_212_FL_OVERRIDE")

for ($i = 0; $i -lt 3; $i++) {
    $r = 61 + $i
    $ws.Cells.Item($r, 1).Value = "CREATE/MODIFY"
    $ws.Cells.Item($r, 2).Value = "LIB_EWS_IT"
    $ws.Cells.Item($r, 3).Value = $newKeys[$i]
    $ws.Cells.Item($r, 5).Value = "String"
    $ws.Cells.Item($r, 6).Value = "String"
}

# Highlight the new rows (A:F) with a yellow fill, matching the inserted-row formatting
$ws.Range("A61:F63").Interior.Color = 65535

# Refresh the view to match the post-edit selection/scroll position
$ws.Range("B64").Select()
$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 1
